# fix(import): add antenne column
#
# Adds a new "antenne" column (K) to the import-template sheet, with a
# sample value of "MONTREUIL" in the example data row, and moves the
# active selection to the next empty cell below the new column so the
# template is ready for further data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (row 1) and sample data (row 2) for the new "antenne" column.
$ws.Range("K1").Value = "antenne"
$ws.Range("K2").Value = "MONTREUIL"

# Reflect the new selection position used after adding the column.
[void]$ws.Range("K3").Select()
